$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 9-14: columns B-H were stored as text (inlineStr, e.g. "500.00") and
# are now re-entered as real numbers so Excel treats them as numeric values.
$numericRows = @{
    9  = @(500, 35500, 5000, 45000, 830, 9500, 126.76)
    10 = @(5000, 40500, 6000, 51000, 840, 10500, 125.93)
    11 = @(5000, 45500, 6000, 57000, 850, 11500, 125.27)
    12 = @(5000, 50500, 6000, 63000, 860, 12500, 124.75)
    13 = @(5000, 55500, 6000, 69000, 870, 13500, 124.32)
    14 = @(5000, 60500, 6000, 75000, 880, 14500, 123.97)
}

$cols = @("B", "C", "D", "E", "F", "G", "H")

foreach ($r in $numericRows.Keys) {
    $values = $numericRows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "{0}{1}" -f $cols[$i], $r
        $ws.Range($addr).Value = $values[$i]
    }
}

# New row 15 with a freshly entered record. These arrive as raw text
# (same shape as the un-normalized rows before this edit), so force the
# cells to text format before assigning so they are kept as strings
# instead of being auto-parsed into numbers.
$ws.Range("A15").Value = "50/84/9000"

$newRowTextCells = @{
    "B15" = "7847.00"
    "C15" = "68347.00"
    "D15" = "7417.00"
    "E15" = "82417.00"
    "F15" = "897.0"
    "G15" = "14070.00"
    "H15" = "120.59"
}

foreach ($addr in $newRowTextCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $newRowTextCells[$addr]
}
